$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the new row style (wrap text, default font/border) to the whole new row first.
$ws.Range("A2:M2").WrapText = $true

# Fill in the new wild-report data row (mirrors the header row's columns).
$ws.Range("A2").Value = "2023/24"
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = "Registered Trapline"
$ws.Range("D2").Value = "TR0515T003"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "No"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""

# G2 (WMU) must stay text ("515"), not be auto-converted to a number. Stage it on a
# scratch cell formatted as text, then bring only the value across so the destination
# keeps the wrap-text style already applied to the rest of the row.
$ws.Range("G99").NumberFormat = "@"
$ws.Range("G99").Value = "515"
$ws.Range("G99").Copy()
$ws.Range("G2").PasteSpecial(-4163, $null)
$ws.Range("G99").Delete()
